$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E stay as TEXT (many values look numeric, e.g. "309.90",
# "0.3662", "26.961.31" -- Excel would otherwise coerce/reformat them).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.961.31"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").Value = "1.818.42"
$ws.Range("E3").Value = "  +0.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "309.90"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").Value = "0.3662"
$ws.Range("E8").Value = "  -1.06%  "

# Row 9
$ws.Range("D9").Value = "0.07354"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").Value = "0.8717"
$ws.Range("E10").Value = "  -0.48%  "

# Row 11
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -0.99%  "

# Row 12
$ws.Range("D12").Value = "1.828.67"
$ws.Range("E12").Value = "  -0.97%  "

# Row 13
$ws.Range("D13").Value = "5.397"
$ws.Range("E13").Value = "  +0.78%  "

# Row 14
$ws.Range("D14").Value = "0.07111"

# Row 15
$ws.Range("D15").Value = "6.506"
$ws.Range("E15").Value = "  -0.03%  "

# Row 16
$ws.Range("D16").Value = "91.45"
$ws.Range("E16").Value = "  -0.22%  "

# Row 17
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").Value = "0.000008695"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").Value = "14.65"

# Row 21
$ws.Range("D21").Value = "26.988.16"

# Row 22
$ws.Range("D22").Value = "5.296"

# Row 23
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "2.047.76"
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").Value = "1.894"
$ws.Range("E25").Value = "  -0.35%  "

# Row 26
$ws.Range("D26").Value = "150.74"
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").Value = "18.44"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28
$ws.Range("D28").Value = "2.138"
$ws.Range("E28").Value = "  -0.64%  "

# Row 29
$ws.Range("D29").Value = "5.240"
$ws.Range("E29").Value = "  -1.69%  "

# Row 30
$ws.Range("D30").Value = "116.48"
$ws.Range("E30").Value = "  +0.65%  "

# Row 31
$ws.Range("D31").Value = "0.08899"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").Value = "0.7587"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("D34").Value = "4.504"
$ws.Range("E34").Value = "  +1.11%  "

# Row 35
$ws.Range("D35").Value = "2.906"
$ws.Range("E35").Value = "  -0.50%  "

# Row 36
$ws.Range("E36").Value = "  +0.10%  "

# Row 37
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("D38").Value = "0.05296"
$ws.Range("E38").Value = "  +1.00%  "

# Row 39
$ws.Range("D39").Value = "0.01947"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40
$ws.Range("D40").Value = "2.969"
$ws.Range("E40").Value = "  +1.42%  "

# Row 41
$ws.Range("D41").Value = "7.172"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("D42").Value = "0.5280"
$ws.Range("E42").Value = "  -0.68%  "

# Row 43
$ws.Range("D43").Value = "2.353"
$ws.Range("E43").Value = "  -3.45%  "

# Row 44
$ws.Range("D44").Value = "0.1658"
$ws.Range("E44").Value = "  -0.32%  "

# Row 45
$ws.Range("D45").Value = "8.446"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("D46").Value = "0.4863"
$ws.Range("E46").Value = "  -2.22%  "

# Row 47
$ws.Range("D47").Value = "10.46"
$ws.Range("E47").Value = "  +1.02%  "

# Row 48
$ws.Range("E48").Value = "  +0.14%  "

# Row 49 (was Quant, now NEARProtocol)
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  -0.08%  "

# Row 50 (was NEARProtocol, now Quant)
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "103.46"
$ws.Range("E50").Value = "  -0.43%  "

# Row 51
$ws.Range("D51").Value = "0.06290"
$ws.Range("E51").Value = "  -0.08%  "
